$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Ancienne Bureautique :" link row (row 9) was removed from the sheet.
# Clear its four cell values while leaving the row's formatting (style) intact.
$ws.Range("A9:D9").ClearContents() | Out-Null

# The user's selection ended up on D9 after the edit.
$ws.Range("D9").Select() | Out-Null
